$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Build the new row 10 ("checking checkout") BEFORE touching H3/D10, since
#    the new I10 cell needs to reuse the distinctive "plain" style that
#    currently lives on the clearCartItems cells (H3 / D10).
# ---------------------------------------------------------------------------

# Copy the whole action sequence (login / productCatalogPage / productDetailPage /
# cartCheck / checkout / orderCOD) from row 6, which already has matching styles,
# into row 10 columns D:I.
$ws.Range("D6:I6").Copy($ws.Range("D10"))

# The trailing action cell (I10, "orderCOD") should use the special "plain"
# style that the old clearCartItems cells (H3) used - copy formatting only.
$ws.Range("H3").Copy()
$ws.Range("I10").PasteSpecial(-4122)

# Fill in the scenario/test-id/flag columns for the new row.
$ws.Range("A10").Value = "checking checkout"
$ws.Range("B10").Value = [string][char]8220 + [string][char]8221
$ws.Range("C10").Value = "YES"

# Final action in the new sequence.
$ws.Range("J10").Value = "confirmationPage"

# ---------------------------------------------------------------------------
# 2. Rename the clearCartItems action to cartCheck (row 3), matching the
#    plain "productDetailPage" style used by its row neighbours.
# ---------------------------------------------------------------------------
$ws.Range("G3").Copy($ws.Range("H3"))
$ws.Range("H3").Value = "cartCheck"

# Remove the now-unused trailing blank cell in row 3.
$ws.Range("I3").Clear()

# ---------------------------------------------------------------------------
# 3. Update the Multi line items test data value.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = [char]8220 + " 431 53" + [char]8221

# ---------------------------------------------------------------------------
# 4. Flip Execution Flag from YES to NO for the already-covered scenarios.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "NO"
$ws.Range("C4").Value = "NO"
$ws.Range("C6").Value = "NO"
$ws.Range("C7").Value = "NO"
$ws.Range("C8").Value = "NO"
